$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new "Winnaar" column, matching the style of the other headers
$ws.Range("G1").Value = "Winnaar"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# New "Winnaar" column values (row 2 through row 20)
$winnaarValues = @(0, 2, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $winnaarValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $winnaarValues[$i]
}
